# Apply "Add working set of sequences" edit:
# Replace the 192 German verb words in column B (rows 2-193) with the new
# working-set word list, in order, matching the target sharedStrings table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'flüchten',
    'bitten',
    'wechseln',
    'münzen',
    'sprengen',
    'saufen',
    'tropfen',
    'pfeifen',
    'knien',
    'liegen',
    'stören',
    'machen',
    'warnen',
    'lenken',
    'sinken',
    'greifen',
    'rasen',
    'bergen',
    'seufzen',
    'hindern',
    'dringen',
    'pflanzen',
    'tollen',
    'suchen',
    'helfen',
    'schwingen',
    'zeugen',
    'scheiden',
    'bellen',
    'wecken',
    'fischen',
    'schwächen',
    'ächzen',
    'weichen',
    'stürmen',
    'nerven',
    'kehren',
    'wehtun',
    'fangen',
    'boxen',
    'rücken',
    'schlucken',
    'stoppen',
    'planen',
    'liefern',
    'töten',
    'schulden',
    'wehren',
    'kosten',
    'läuten',
    'wenden',
    'brauchen',
    'sterben',
    'spinnen',
    'nennen',
    'schreiten',
    'altern',
    'lassen',
    'leisten',
    'drehen',
    'fassen',
    'pflegen',
    'äußern',
    'fließen',
    'zeigen',
    'wahren',
    'runden',
    'klettern',
    'loben',
    'ehren',
    'heilen',
    'treiben',
    'laufen',
    'sorgen',
    'danken',
    'kichern',
    'quälen',
    'klappen',
    'sperren',
    'graben',
    'arten',
    'knarren',
    'posten',
    'freuen',
    'ändern',
    'rufen',
    'spüren',
    'süßen',
    'trauen',
    'zünden',
    'lesen',
    'schnellen',
    'stärken',
    'deuten',
    'heben',
    'stellen',
    'jagen',
    'leeren',
    'formen',
    'streichen',
    'segnen',
    'mauern',
    'bauen',
    'klagen',
    'folgen',
    'bluten',
    'regeln',
    'wundern',
    'führen',
    'ärgern',
    'bleiben',
    'gelten',
    'stehlen',
    'betteln',
    'schneiden',
    'trennen',
    'decken',
    'schlagen',
    'klingen',
    'reizen',
    'lockern',
    'grüßen',
    'schleppen',
    'zielen',
    'wirken',
    'flehen',
    'werfen',
    'scheitern',
    'trösten',
    'gründen',
    'heulen',
    'stecken',
    'spielen',
    'dürfen',
    'werden',
    'handeln',
    'schmecken',
    'sagen',
    'hören',
    'filmen',
    'wetten',
    'räumen',
    'feiern',
    'weigern',
    'siegen',
    'schämen',
    'träumen',
    'schwören',
    'erben',
    'schauen',
    'platzen',
    'biegen',
    'jubeln',
    'klingeln',
    'lügen',
    'wüten',
    'doppeln',
    'fragen',
    'schalten',
    'scheinen',
    'sichern',
    'zögern',
    'tauchen',
    'stammen',
    'mögen',
    'meistern',
    'leihen',
    'narren',
    'wachsen',
    'hauen',
    'schaden',
    'fallen',
    'enden',
    'stillen',
    'löschen',
    'schenken',
    'dienen',
    'schrecken',
    'geben',
    'mühen',
    'irren',
    'buchen',
    'fällen',
    'füttern',
    'backen',
    'achten',
    'foltern',
    'fahren',
    'kümmern',
    'malen',
    'zahlen',
    'kürzen'
)
for ($i = 0; $i -lt $words.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $words[$i]
}
